$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----- Header row (row 73): Task Type | Count | Total Time | Waiting Time -----
$ws.Range("E73").Value = "Task Type"
$ws.Range("F73").Value = "Count"
$ws.Range("G73").Value = "Total Time"
$ws.Range("H73").Value = "Waiting Time"

# ----- Transportation summary row (row 74) -----
$ws.Range("E74").Value = "Transport"
$ws.Range("F74").Formula = '=COUNTIF(A2:A62,"transportation")'
$ws.Range("G74").Formula = '=AVERAGEIF(A2:A62,"transportation",F2:F62)'
$ws.Range("H74").Formula = '=AVERAGEIF(A2:A62,"transportation",H2:H62)'

# ----- Charging summary row (row 75) -----
$ws.Range("E75").Value = "Charging"
$ws.Range("F75").Formula = '=COUNTIF(A2:A62,"charging")'
$ws.Range("G75").Formula = '=AVERAGEIF(A2:A62,"charging",F2:F62)'
$ws.Range("H75").Formula = '=AVERAGEIF(A2:A62,"charging",H2:H62)'

# ----- Formatting -----
# Whole block: thin box border around every cell + centered text
$header = $ws.Range("E73:H73")
$body = $ws.Range("E74:H75")
$block = $ws.Range("E73:H75")
$block.Borders.LineStyle = 1
$block.HorizontalAlignment = -4108

# Header row gets extra: bold font + light-blue ("Accent1, Lighter 40%") fill
$header.Font.Bold = $true
$header.Interior.Color = 15652797

# ----- View niceties matching the target selection/scroll -----
$ws.Range("D79").Select()

Write-Host "edit applied"
